$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 17859504
$ws.Range("I62").Value = 25001846
$ws.Range("K62").Value = 25001846
$ws.Range("M62").Value = -25001222
$ws.Range("H65").Value = 17859504
$ws.Range("I65").Value = 25001846
$ws.Range("K65").Value = 125009230
$ws.Range("M65").Value = -125006110
$ws.Range("H74").Value = 11993.692
$ws.Range("I74").Value = 10909.833
$ws.Range("K74").Value = 10909.833
$ws.Range("M74").Value = -9973.833000000001
$ws.Range("H77").Value = 11993.692
$ws.Range("I77").Value = 10909.833
$ws.Range("K77").Value = 54549.165
$ws.Range("M77").Value = -49869.165
$ws.Range("H125").Value = 8097.143
$ws.Range("I125").Value = 7507.5
$ws.Range("K125").Value = 67567.5
$ws.Range("M125").Value = -65107.5
$ws.Range("H132").Value = 4686.943
$ws.Range("I132").Value = 4765.893
$ws.Range("K132").Value = 14297.679
$ws.Range("M132").Value = -11767.679

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2729.3674
$ws.Range("I32").Value = 2488.0852
$ws.Range("K32").Value = 2488.0852
$ws.Range("M32").Value = -2201.0852
$ws.Range("H102").Value = 5266.4287
$ws.Range("I102").Value = 1717
$ws.Range("K102").Value = 1717
$ws.Range("M102").Value = -95
$ws.Range("H132").Value = 2504.1667
$ws.Range("I132").Value = 848.5
$ws.Range("K132").Value = 2545.5
$ws.Range("M132").Value = -15.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 45224
$ws.Range("I105").Value = 100957.8
$ws.Range("K105").Value = 100957.8
$ws.Range("M105").Value = -99210.8

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 30000400
$ws.Range("I4").Value = 5000600
$ws.Range("K4").Value = 5000600
$ws.Range("M4").Value = -5000488
$ws.Range("H6").Value = 3578.6
$ws.Range("I6").Value = 3578.6
$ws.Range("K6").Value = 3578.6
$ws.Range("M6").Value = -3465.6
$ws.Range("H31").Value = 3481.4375
$ws.Range("I31").Value = 2055.5715
$ws.Range("K31").Value = 2055.5715
$ws.Range("M31").Value = -1760.5715
$ws.Range("H34").Value = 3481.4375
$ws.Range("I34").Value = 2055.5715
$ws.Range("K34").Value = 2055.5715
$ws.Range("M34").Value = -1853.5715
$ws.Range("H60").Value = 64000
$ws.Range("J60").Value = 64000
$ws.Range("L60").Value = 64000
$ws.Range("N60").Value = -65022
$ws.Range("H96").Value = 39906
$ws.Range("J96").Value = 39906
$ws.Range("L96").Value = 39906
$ws.Range("N96").Value = -45398
$ws.Range("H134").Value = 4647.8887
$ws.Range("I134").Value = 3835.2856
$ws.Range("K134").Value = 11505.8568
$ws.Range("M134").Value = -8970.856800000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23848098
$ws.Range("I4").Value = 25594430
$ws.Range("K4").Value = 76783290
$ws.Range("M4").Value = -76783178
$ws.Range("H86").Value = 1375482.2
$ws.Range("I86").Value = 634
$ws.Range("J86").Value = 2200391.2
$ws.Range("K86").Value = 1902
$ws.Range("L86").Value = 6601173.600000001
$ws.Range("M86").Value = -716
$ws.Range("N86").Value = -6603545.600000001
$ws.Range("H89").Value = 1375482.2
$ws.Range("I89").Value = 634
$ws.Range("J89").Value = 2200391.2
$ws.Range("K89").Value = 5706
$ws.Range("L89").Value = 19803520.8
$ws.Range("M89").Value = 222
$ws.Range("N89").Value = -19815376.8
$ws.Range("H107").Value = 128777
$ws.Range("J107").Value = 137288.14
$ws.Range("L107").Value = 411864.42
$ws.Range("N107").Value = -415704.42
$ws.Range("H131").Value = 3016.718
$ws.Range("J131").Value = 3846.7693
$ws.Range("L131").Value = 11540.3079
$ws.Range("N131").Value = -21620.3079
$ws.Range("H132").Value = 1864.862
$ws.Range("J132").Value = 4812.25
$ws.Range("L132").Value = 43310.25
$ws.Range("N132").Value = -48370.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1945.1111
$ws.Range("I13").Value = 2400
$ws.Range("K13").Value = 2400
$ws.Range("M13").Value = -2261
$ws.Range("H122").Value = 82121.69500000001
$ws.Range("I122").Value = 82121.69500000001
$ws.Range("K122").Value = 246365.085
$ws.Range("M122").Value = -243915.085
$ws.Range("H132").Value = 219948.98
$ws.Range("I132").Value = 272520.3
$ws.Range("J132").Value = 3822.3333
$ws.Range("K132").Value = 817560.8999999999
$ws.Range("L132").Value = 11466.9999
$ws.Range("M132").Value = -815030.8999999999
$ws.Range("N132").Value = -16526.9999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 422992.03
$ws.Range("I7").Value = 672560.8
$ws.Range("K7").Value = 672560.8
$ws.Range("M7").Value = -672448.8
$ws.Range("H22").Value = 708.7143
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1590
$ws.Range("H27").Value = 708.7143
$ws.Range("J27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("N27").Value = -1214
$ws.Range("H46").Value = 3303.25
$ws.Range("I46").Value = 2352.0588
$ws.Range("K46").Value = 2352.0588
$ws.Range("M46").Value = -2164.0588
$ws.Range("H55").Value = 474.94446
$ws.Range("I55").Value = 299.34616
$ws.Range("J55").Value = 931.5
$ws.Range("K55").Value = 299.34616
$ws.Range("L55").Value = 931.5
$ws.Range("M55").Value = -126.34616
$ws.Range("N55").Value = -1277.5
$ws.Range("H68").Value = 136482.5
$ws.Range("I68").Value = 4974.2
$ws.Range("K68").Value = 4974.2
$ws.Range("M68").Value = -4225.2
$ws.Range("H71").Value = 136482.5
$ws.Range("I71").Value = 4974.2
$ws.Range("K71").Value = 24871
$ws.Range("M71").Value = -21127
$ws.Range("H122").Value = 557645.9
$ws.Range("I122").Value = 404407.44
$ws.Range("J122").Value = 852335.1
$ws.Range("K122").Value = 1213222.32
$ws.Range("L122").Value = 2557005.3
$ws.Range("M122").Value = -1210772.32
$ws.Range("N122").Value = -2561905.3
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 422992.03
$ws.Range("I126").Value = 672560.8
$ws.Range("K126").Value = 2017682.4
$ws.Range("M126").Value = -2015212.4
$ws.Range("H136").Value = 4638.9546
$ws.Range("I136").Value = 3009.2222
$ws.Range("K136").Value = 9027.6666
$ws.Range("M136").Value = -6477.6666

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 2500
$ws.Range("I49").Value = 2500
$ws.Range("K49").Value = 2500
$ws.Range("M49").Value = -2270
$ws.Range("H107").Value = 73969.36
$ws.Range("I107").Value = 85939.25
$ws.Range("J107").Value = 2150
$ws.Range("K107").Value = 257817.75
$ws.Range("L107").Value = 6450
$ws.Range("M107").Value = -255897.75
$ws.Range("N107").Value = -10290
$ws.Range("H122").Value = 23812952
$ws.Range("I122").Value = 35716696
$ws.Range("K122").Value = 107150088
$ws.Range("M122").Value = -107147638
$ws.Range("H133").Value = 64996
$ws.Range("J133").Value = 64996
$ws.Range("L133").Value = 64996
$ws.Range("N133").Value = -75116
$ws.Range("H136").Value = 336417.12
$ws.Range("J136").Value = 4930.625
$ws.Range("L136").Value = 14791.875
$ws.Range("N136").Value = -19891.875
